$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48875.19
$ws.Range("I2").Value = 1335.7368
$ws.Range("K2").Value = 1335.7368
$ws.Range("M2").Value = -1222.7368

$ws.Range("H32").Value = 6723.69
$ws.Range("I32").Value = 5370.7803
$ws.Range("J32").Value = 20403.111
$ws.Range("K32").Value = 5370.7803
$ws.Range("L32").Value = 20403.111
$ws.Range("M32").Value = -5083.7803
$ws.Range("N32").Value = -20977.111

$ws.Range("H97").Value = 33949.71
$ws.Range("I97").Value = 48750
$ws.Range("J97").Value = 2869.1
$ws.Range("K97").Value = 48750
$ws.Range("L97").Value = 2869.1
$ws.Range("M97").Value = -48254
$ws.Range("N97").Value = -3861.1

$ws.Range("H110").Value = 77085340
$ws.Range("I110").Value = 77085340
$ws.Range("K110").Value = 77085340
$ws.Range("M110").Value = -77083295

$ws.Range("H116").Value = 48875.19
$ws.Range("I116").Value = 1335.7368
$ws.Range("K116").Value = 1335.7368
$ws.Range("M116").Value = 958.2632000000001

$ws.Range("H122").Value = 1926.0416
$ws.Range("I122").Value = 1712.25
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 5136.75
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -2686.75
$ws.Range("N122").Value = -13885

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48875.19
$ws.Range("I3").Value = 1335.7368
$ws.Range("K3").Value = 1335.7368
$ws.Range("M3").Value = -1221.7368

$ws.Range("H25").Value = 9700
$ws.Range("I25").Value = 2266.6667
$ws.Range("J25").Value = 32000
$ws.Range("K25").Value = 2266.6667
$ws.Range("L25").Value = 32000
$ws.Range("M25").Value = -2031.6667
$ws.Range("N25").Value = -32470

$ws.Range("H105").Value = 144923.64
$ws.Range("I105").Value = 113083.336
$ws.Range("J105").Value = 202236.2
$ws.Range("K105").Value = 113083.336
$ws.Range("L105").Value = 202236.2
$ws.Range("M105").Value = -111336.336
$ws.Range("N105").Value = -205730.2

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 299.2
$ws.Range("I22").Value = 186.5
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 186.5
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = 163.5
$ws.Range("N22").Value = -1450

$ws.Range("H58").Value = 1838.2188
$ws.Range("I58").Value = 1620.579
$ws.Range("J58").Value = 2156.3076
$ws.Range("K58").Value = 1620.579
$ws.Range("L58").Value = 2156.3076
$ws.Range("M58").Value = -1417.579
$ws.Range("N58").Value = -2562.3076

$ws.Range("H99").Value = 2829.7334
$ws.Range("I99").Value = 2898.2
$ws.Range("J99").Value = 2795.5
$ws.Range("K99").Value = 2898.2
$ws.Range("L99").Value = 2795.5
$ws.Range("M99").Value = -1400.2
$ws.Range("N99").Value = -5791.5

$ws.Range("H122").Value = 1174.5714
$ws.Range("I122").Value = 1164.4
$ws.Range("K122").Value = 3493.2
$ws.Range("M122").Value = -1043.2

$ws.Range("H126").Value = 2829.7334
$ws.Range("I126").Value = 2898.2
$ws.Range("J126").Value = 2795.5
$ws.Range("K126").Value = 8694.599999999999
$ws.Range("L126").Value = 8386.5
$ws.Range("M126").Value = -6224.599999999999
$ws.Range("N126").Value = -13326.5

$ws.Range("H132").Value = 93754190
$ws.Range("I132").Value = 125005780
$ws.Range("J132").Value = 62502616
$ws.Range("K132").Value = 375017340
$ws.Range("L132").Value = 187507848
$ws.Range("M132").Value = -375014810
$ws.Range("N132").Value = -187512908

$ws.Range("H136").Value = 1838.2188
$ws.Range("I136").Value = 1620.579
$ws.Range("J136").Value = 2156.3076
$ws.Range("K136").Value = 4861.737
$ws.Range("L136").Value = 6468.9228
$ws.Range("M136").Value = -2311.737
$ws.Range("N136").Value = -11568.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 746.25
$ws.Range("J61").Value = 746.25
$ws.Range("L61").Value = 2238.75
$ws.Range("N61").Value = -2668.75

$ws.Range("H98").Value = 72562.71000000001
$ws.Range("I98").Value = 420.5
$ws.Range("J98").Value = 101419.6
$ws.Range("K98").Value = 1261.5
$ws.Range("L98").Value = 304258.8
$ws.Range("M98").Value = 236.5
$ws.Range("N98").Value = -307254.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 67891.84
$ws.Range("I70").Value = 95436.45
$ws.Range("J70").Value = 7293.7
$ws.Range("K70").Value = 95436.45
$ws.Range("L70").Value = 7293.7
$ws.Range("M70").Value = -95166.45
$ws.Range("N70").Value = -7833.7

$ws.Range("H73").Value = 67891.84
$ws.Range("I73").Value = 95436.45
$ws.Range("J73").Value = 7293.7
$ws.Range("K73").Value = 95436.45
$ws.Range("L73").Value = 7293.7
$ws.Range("M73").Value = -94500.45
$ws.Range("N73").Value = -9165.700000000001

$ws.Range("H80").Value = 91005330
$ws.Range("I80").Value = 166840000
$ws.Range("J80").Value = 3733.2
$ws.Range("K80").Value = 166840000
$ws.Range("L80").Value = 3733.2
$ws.Range("M80").Value = -166839002
$ws.Range("N80").Value = -5729.2

$ws.Range("H83").Value = 91005330
$ws.Range("I83").Value = 166840000
$ws.Range("J83").Value = 3733.2
$ws.Range("K83").Value = 834200000
$ws.Range("L83").Value = 18666
$ws.Range("M83").Value = -834195008
$ws.Range("N83").Value = -28650

$ws.Range("H102").Value = 3269.238
$ws.Range("I102").Value = 1974.5454
$ws.Range("J102").Value = 4693.4
$ws.Range("K102").Value = 1974.5454
$ws.Range("L102").Value = 4693.4
$ws.Range("M102").Value = -352.5454
$ws.Range("N102").Value = -7937.4

$ws.Range("H113").Value = 1911.091
$ws.Range("I113").Value = 1860.2858
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1860.2858
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 309.7141999999999
$ws.Range("N113").Value = -6340

$ws.Range("H122").Value = 2841.0833
$ws.Range("I122").Value = 2286.3333
$ws.Range("J122").Value = 4505.3335
$ws.Range("K122").Value = 6858.999899999999
$ws.Range("L122").Value = 13516.0005
$ws.Range("M122").Value = -4408.999899999999
$ws.Range("N122").Value = -18416.0005

$ws.Range("H136").Value = 16251.842
$ws.Range("J136").Value = 16251.842
$ws.Range("L136").Value = 48755.526
$ws.Range("N136").Value = -53855.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1892.381
$ws.Range("I61").Value = 1898.0769
$ws.Range("J61").Value = 1883.125
$ws.Range("K61").Value = 1898.0769
$ws.Range("L61").Value = 1883.125
$ws.Range("M61").Value = -1696.0769
$ws.Range("N61").Value = -2287.125

$ws.Range("H113").Value = 1892.381
$ws.Range("I113").Value = 1898.0769
$ws.Range("J113").Value = 1883.125
$ws.Range("K113").Value = 1898.0769
$ws.Range("L113").Value = 1883.125
$ws.Range("M113").Value = 271.9231
$ws.Range("N113").Value = -6223.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1975.0294
$ws.Range("I122").Value = 1505.4615
$ws.Range("J122").Value = 3501.125
$ws.Range("K122").Value = 4516.3845
$ws.Range("L122").Value = 10503.375
$ws.Range("M122").Value = -2066.3845
$ws.Range("N122").Value = -15403.375

$ws.Range("H132").Value = 4395.136
$ws.Range("I132").Value = 5276.0713
$ws.Range("J132").Value = 2853.5
$ws.Range("K132").Value = 15828.2139
$ws.Range("L132").Value = 8560.5
$ws.Range("M132").Value = -13298.2139
$ws.Range("N132").Value = -13620.5

$ws.Range("H136").Value = 2225.5938
$ws.Range("I136").Value = 953.6
$ws.Range("J136").Value = 3347.9412
$ws.Range("K136").Value = 2860.8
$ws.Range("L136").Value = 10043.8236
$ws.Range("M136").Value = -310.8000000000002
$ws.Range("N136").Value = -15143.8236
